# [Kadastro App] Yeni kayit eklendi: 3009
# Appends the new record (Kayit No 3009) as the next row (68) on both the
# master "Kayitlar" sheet and the per-birim "Erdemli" sheet (the record's
# Birim). Columns A/B/D hold digit-only-looking text ("3009", date string,
# parcel count) that Excel would otherwise auto-coerce to Number/Date, so
# they are entered with a leading apostrophe (exactly like a user typing
# '3009 into a General-formatted cell) to keep them stored as Text, matching
# the rest of the sheet (which relies on the numberStoredAsText ignored
# error to suppress the green-triangle warning).

$wb = $excel.ActiveWorkbook

$newRow = 68
$kayitNo   = "3009"
$tarih     = "2025-09-11"
$birim     = "Erdemli"
$parsel    = "1"
$is        = "CİNS DEĞ."
$personel  = "AYHAN KARADAYI (K.Teknisyeni), EMİNE ALANLI KIRCILI (K.Mühendisi)"

foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = "'" + $kayitNo
    $ws.Cells.Item($newRow, 2).Value = "'" + $tarih
    $ws.Cells.Item($newRow, 3).Value = $birim
    $ws.Cells.Item($newRow, 4).Value = "'" + $parsel
    $ws.Cells.Item($newRow, 5).Value = $is
    $ws.Cells.Item($newRow, 6).Value = $personel
}
